$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hide future modules: set link_it (column C) to FALSE for rows 3-20,
# except row 5 which is already FALSE.
$rows = 3,4,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20
foreach ($r in $rows) {
    $ws.Range("C$r").Value = $false
}

# Update the active selection to E7
$ws.Range("E7").Select()
